$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 92

$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025-10-21"

$cellB = $ws.Cells.Item($row, 2)
$cellB.NumberFormat = "@"
$cellB.Value = "21:22:08"

$cellC = $ws.Cells.Item($row, 3)
$cellC.Value = "1.00 EUR = 1,806.1939"
